$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (this also updates the _xlnm._FilterDatabase defined name
# and the workbook's sheet list automatically).
$ws.Name = "CDM Data Catalogue"

# Update the Barrage Flow record's Document Title and Date to reflect the new
# source file / extended date range.
$ws.Range("F18").Value = "barrage_daily_total.csv"
$ws.Range("I18").Value = "1990 - 2022"

# Restore the active cell selection to reflect where the author left off.
$ws.Range("K20").Select()
